# Weekly update: a new "Femacal de La Calera - Espárragos" record needs to be
# inserted right after the row dated 2022-12-26 (serial 44900, row 14), which
# pushes every subsequent weekly record down by one row (old row 59 becomes
# row 60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 15, shifting rows 15:59 down to 16:60.
$ws.Range("A15:R15").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44910
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Verde"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 650
$ws.Range("K15").Value = 1500
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1500
$ws.Range("N15").Value = "$/kilo"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 1500
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"
